# fix profile pictures (in Process)
# Rename the "Rooms" worksheet to "Profile Picture Map"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rooms")
$ws.Name = "Profile Picture Map"
